$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/no-format text cell) used to strip the transient
# "@" number-format style that gets attached when we force text-typed values
# below, so the edited cells keep the same style index as in the source file.
$textStyle = $ws.Range("B2").Style

function Set-TextValue([string]$cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $textStyle
}

$ws.Range("D2").Value = "67.680.09"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "3.505.11"
$ws.Range("E3").Value = "  -0.52%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue "D5" "607.56"
$ws.Range("E5").Value = "  -0.98%  "
Set-TextValue "D6" "152.04"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "3.503.83"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue "D9" "0.485"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +2.86%  "
Set-TextValue "D11" "7.59"
$ws.Range("E11").Value = "  +6.79%  "
$ws.Range("E12").Value = "  +1.29%  "
Set-TextValue "D13" "0.0000218"
$ws.Range("E13").Value = "  -1.49%  "
Set-TextValue "D14" "32.37"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "4.093.73"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "3.504.51"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "67.543.65"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  -0.62%  "
Set-TextValue "D19" "6.52"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("E20").Value = "  +1.92%  "
Set-TextValue "D21" "9.89"
$ws.Range("E21").Value = "  +4.64%  "
Set-TextValue "D22" "447.40"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +0.65%  "
Set-TextValue "D24" "78.13"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "3.640.16"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -0.01%  "
Set-TextValue "D27" "0.0000126"
$ws.Range("E27").Value = "  -2.96%  "
Set-TextValue "D28" "10.11"
$ws.Range("E28").Value = "  -2.21%  "
Set-TextValue "D29" "8.75"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("E30").Value = "  +0.39%  "
Set-TextValue "D31" "1.64"
$ws.Range("E31").Value = "  +6.71%  "
Set-TextValue "D32" "0.169"
$ws.Range("E32").Value = "  +3.86%  "
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  -0.10%  "
Set-TextValue "D34" "25.64"
$ws.Range("E34").Value = "  -1.06%  "
Set-TextValue "D35" "6.17"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").Value = "3.493.69"
$ws.Range("E37").Value = "  -0.59%  "
Set-TextValue "D38" "7.99"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +7.07%  "
Set-TextValue "D41" "179.13"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  -0.08%  "
Set-TextValue "D43" "0.0893"
$ws.Range("E43").Value = "  +1.40%  "
Set-TextValue "D44" "5.47"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "30.32"
$ws.Range("E45").Value = "  +5.94%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D46" "0.889"
$ws.Range("E46").Value = "  +0.80%  "
Set-TextValue "D47" "46.29"
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("E48").Value = "  +2.73%  "
Set-TextValue "D49" "2.57"
$ws.Range("E49").Value = "  -1.92%  "
Set-TextValue "D50" "7.63"
$ws.Range("E50").Value = "  +0.53%  "
Set-TextValue "D51" "0.993"
$ws.Range("E51").Value = "  -0.45%  "
